# feat: add 2022-Q1 data
#
# Before:  Sheet1 = "2021-Q4" (fund holdings), Sheet2 = "总计" (totals: 2021-Q4 only)
# After:   Sheet1 = "2021-Q4" (unchanged), Sheet2 = "2022-Q1" (new fund holdings),
#          Sheet3 = "总计" (totals: 2022-Q1 + 2021-Q4)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "2021-Q4" - stays as-is
$wsOldTotal = $wb.Worksheets.Item(2)   # currently "总计"

# --- Step 1: free up the "总计" name by renaming this sheet to "2022-Q1" ---
$wsOldTotal.Name = "2022-Q1"
$wsQ1 = $wsOldTotal

# --- Step 2: insert the new "总计" sheet right after "2022-Q1" ---
$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

# Match the page-margin convention used by the rest of the workbook
# (0.75in/0.75in/1in/1in/0.5in/0.5in == 54/54/72/72/36/36 points).
$wsTotal.PageSetup.LeftMargin = 54
$wsTotal.PageSetup.RightMargin = 54
$wsTotal.PageSetup.TopMargin = 72
$wsTotal.PageSetup.BottomMargin = 72
$wsTotal.PageSetup.HeaderMargin = 36
$wsTotal.PageSetup.FooterMargin = 36

# =========================================================================
# Build the "2022-Q1" fund-holdings sheet (same layout/style as "2021-Q4")
# =========================================================================
$wsQ1.Cells.Clear()

# Header row (copy style+text from the "2021-Q4" header so formatting matches)
$ws1.Range("B1:H1").Copy($wsQ1.Range("B1:H1"))
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# Column A (index numbers) - copy the bordered/bold style down for all 14 rows
$ws1.Range("A2").Copy($wsQ1.Range("A2:A15"))

# Columns B-G must stay as literal text (even the numeric-looking ones), so
# force a Text number format before assigning their values.
$wsQ1.Range("B2:G15").NumberFormat = "@"

$data = @(
    @("007689","国投瑞银新能源混合A","78.74","91.13","4.10","3.2283",10),
    @("012148","国投瑞银产业趋势混合型证券投资基金A","45.53","92.28","3.70","1.6846",9),
    @("012079","信达澳银新能源精选混合型证券投资基金","35.11","94.31","4.46","1.5659",10),
    @("007690","国投瑞银新能源混合C","37.84","91.13","4.10","1.5514",10),
    @("001704","国投瑞银进宝灵活配置混合","33.25","92.49","4.31","1.4331",7),
    @("310328","申万菱信新动力混合","34.56","73.99","2.82","0.9746",4),
    @("012149","国投瑞银产业趋势混合型证券投资基金C","18.42","92.28","3.70","0.6815",9),
    @("013513","长安先进制造混合A","2.95","91.34","3.68","0.1086",8),
    @("002293","南方益和灵活配置混合","1.40","83.59","3.55","0.0497",8),
    @("013514","长安先进制造混合C","0.52","91.34","3.68","0.0191",8),
    @("001261","中融新机遇灵活配置混合","0.34","93.08","5.22","0.0177",10),
    @("004557","北信瑞丰鼎丰灵活配置混合","0.39","64.13","4.53","0.0177",9),
    @("011800","申万菱信价值精选混合型证券投资基金","0.57","81.46","3.05","0.0174",5),
    @("005536","渤海汇金量化成长混合","0.61","88.57","1.07","0.0065",3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $wsQ1.Cells.Item($row, 1).Value = $i
    $wsQ1.Cells.Item($row, 2).Value = $rec[0]
    $wsQ1.Cells.Item($row, 3).Value = $rec[1]
    $wsQ1.Cells.Item($row, 4).Value = $rec[2]
    $wsQ1.Cells.Item($row, 5).Value = $rec[3]
    $wsQ1.Cells.Item($row, 6).Value = $rec[4]
    $wsQ1.Cells.Item($row, 7).Value = $rec[5]
    $wsQ1.Cells.Item($row, 8).Value = $rec[6]
}

# =========================================================================
# Build the "总计" (totals) sheet
# =========================================================================
$wsTotal.Cells.Clear()

$ws1.Range("B1:D1").Copy($wsTotal.Range("B1:D1"))
$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$ws1.Range("A2").Copy($wsTotal.Range("A2:A3"))
$wsTotal.Range("B2:B3").NumberFormat = "@"

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 14
$wsTotal.Cells.Item(2, 4).Value = 11.36

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2021-Q4"
$wsTotal.Cells.Item(3, 3).Value = 12
$wsTotal.Cells.Item(3, 4).Value = 4.42

# Restore the original active sheet/selection ("2021-Q4" stays selected).
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
